$d = $word.ActiveDocument

# Remove the existing "_GoBack" bookmark so we can re-create it, collapsed,
# at the end of the new first paragraph instead of the current paragraph.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# Insert the new paragraph "Hey newbranch" before all existing content. This
# splits the start of the document into its own paragraph while leaving the
# original paragraph (now second) completely untouched.
$startRange = $d.Range(0, 0)
$startRange.InsertBefore("Hey newbranch" + [char]13)

# A temporary marker character is inserted right after the new paragraph's
# text so that the bookmark we add next is not positioned at the very last
# offset of the paragraph (which would otherwise make it snap to wrap the
# whole paragraph). We then delete the marker, leaving the bookmark
# collapsed cleanly at the end of the paragraph's text.
$newParaTextEnd = "Hey newbranch".Length
$markerRange = $d.Range($newParaTextEnd, $newParaTextEnd)
$markerRange.InsertAfter("Z")

$bookmarkRange = $d.Range($newParaTextEnd, $newParaTextEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$markerRange = $d.Range($newParaTextEnd, $newParaTextEnd + 1)
$markerRange.Delete()
